# Update "想去人数" (want-to-go count) figures across sheets, reflecting the
# site's regenerated data snapshot (gh-pages output at 456a3b4).

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (Exhibitions) -------------------------------------------
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F7").Value  = 91
$ws.Range("F8").Value  = 10163
$ws.Range("F9").Value  = 57
$ws.Range("F10").Value = 3506
$ws.Range("F12").Value = 2438
$ws.Range("F13").Value = 28
$ws.Range("F14").Value = 2793
$ws.Range("F17").Value = 2160
$ws.Range("F19").Value = 93
$ws.Range("F21").Value = 384
$ws.Range("F23").Value = 136
$ws.Range("F26").Value = 216
$ws.Range("F28").Value = 1312
$ws.Range("F29").Value = 9
$ws.Range("F30").Value = 1250
$ws.Range("F31").Value = 103
$ws.Range("F34").Value = 3283
$ws.Range("F35").Value = 3078
$ws.Range("F36").Value = 26
$ws.Range("F38").Value = 1035
$ws.Range("F42").Value = 92
$ws.Range("F46").Value = 39
$ws.Range("F47").Value = 7

# --- Sheet "本地生活" (Local life) ------------------------------------------
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 742
$ws.Range("F3").Value = 984
$ws.Range("F5").Value = 1991

# --- Sheet "全部类型" (All types) -------------------------------------------
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value  = 984
$ws.Range("F9").Value  = 91
$ws.Range("F10").Value = 10163
$ws.Range("F11").Value = 57
$ws.Range("F12").Value = 3506
$ws.Range("F14").Value = 28
$ws.Range("F17").Value = 2160
$ws.Range("F19").Value = 93
$ws.Range("F21").Value = 136
$ws.Range("F24").Value = 216
$ws.Range("F26").Value = 1312
$ws.Range("F27").Value = 9
$ws.Range("F28").Value = 1250
$ws.Range("F33").Value = 3283
$ws.Range("F34").Value = 3078
$ws.Range("F35").Value = 26
$ws.Range("F36").Value = 1035
$ws.Range("F44").Value = 92
$ws.Range("F47").Value = 39
$ws.Range("F48").Value = 7
